{"js": "// Applies the diary-TODO edit:\n//  1) Merge the \"Jeg\" + \" har \" runs (and drop the proofErr spell-check\n//     markers around \"Jeg\") into a single run \"Jeg har \" in the\n//     \".env fil der ikke vil loade\" day's \"haft problemer\" paragraph.\n//  2) In the 17/5 TODO list: \"repository\" -> \"database CRUD for user\",\n//     and \"REST api\" -> \"REST \" + \"user API\" (two runs).\n//  3) Append a \". \" run to the \"Fik endelig hul igennem...\" paragraph and\n//     insert two new paragraphs after it (before the next, empty one).\n\nconst NS_PKG = \"http://schemas.microsoft.com/office/2006/xmlPackage\";\nconst NS_W = \"http://schemas.openxmlformats.org/wordprocessingml/2006/main\";\n\n// Wrap a <w:body> fragment into the pkg:package envelope insertOoxml expects.\nfunction wrapOoxml(bodyXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"' + NS_PKG + '\">' +\n      '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n          '<w:document xmlns:w=\"' + NS_W + '\">' +\n            '<w:body>' + bodyXml + '</w:body>' +\n          '</w:document>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// ---- Change 1: \"Jeg\" + \" har \" -> single run \"Jeg har \", proofErr removed ----\nconst jegHarIdx = items.findIndex(\n  (p) => p.text.indexOf(\"haft problemer med at f\\u00e5 lavet noget ordentlig kode\") !== -1\n);\nif (jegHarIdx === -1) {\n  throw new Error(\"Could not locate the 'Jeg har haft problemer...' paragraph\");\n}\n{\n  const p = items[jegHarIdx];\n  const tail = p.text.substring(\"Jeg har \".length);\n  const newParaXml =\n    \"<w:p>\" +\n      '<w:r><w:rPr><w:lang w:val=\"en-DK\"/></w:rPr><w:t xml:space=\"preserve\">Jeg har </w:t></w:r>' +\n      \"<w:r><w:t>\" + tail + \"</w:t></w:r>\" +\n    \"</w:p>\";\n  p.getRange().insertOoxml(wrapOoxml(newParaXml), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Reload paragraphs since the collection/indices can shift after the edit above.\nparagraphs.load(\"text\");\nawait context.sync();\nconst items2 = paragraphs.items;\n\n// ---- Change 2: TODO list \"repository\" / \"REST api\" (17/5 section) ----\nconst repoIdx = items2.findIndex(\n  (p, i) =>\n    p.text === \"repository\" &&\n    i + 1 < items2.length &&\n    items2[i + 1].text === \"REST api\"\n);\nif (repoIdx === -1) {\n  throw new Error(\"Could not locate the 'repository' / 'REST api' TODO pair\");\n}\nconst restIdx = repoIdx + 1;\n\nconst listPPr =\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>';\n\n{\n  const repoP = items2[repoIdx];\n  const newXml =\n    \"<w:p>\" + listPPr + \"<w:r><w:t>database CRUD for user</w:t></w:r></w:p>\";\n  repoP.getRange().insertOoxml(wrapOoxml(newXml), Word.InsertLocation.replace);\n  await context.sync();\n}\n{\n  // Re-fetch the REST api paragraph range fresh (indices unaffected: same count of paragraphs).\n  const restP = paragraphs.items[restIdx];\n  const newXml =\n    \"<w:p>\" +\n      listPPr +\n      '<w:r><w:t xml:space=\"preserve\">REST </w:t></w:r>' +\n      \"<w:r><w:t>user API</w:t></w:r>\" +\n    \"</w:p>\";\n  restP.getRange().insertOoxml(wrapOoxml(newXml), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---- Change 3: append \". \" run + two new paragraphs after \"Fik endelig...\" ----\nparagraphs.load(\"text\");\nawait context.sync();\nconst items3 = paragraphs.items;\nconst fikIdx = items3.findIndex((p) => p.text.indexOf(\"Fik endelig hul igennem\") === 0);\nif (fikIdx === -1) {\n  throw new Error(\"Could not locate the 'Fik endelig hul igennem...' paragraph\");\n}\n{\n  const fikP = items3[fikIdx];\n  const originalText = fikP.text;\n  const newXml =\n    \"<w:p><w:r><w:t>\" + originalText + '</w:t></w:r><w:r><w:t xml:space=\"preserve\">. </w:t></w:r></w:p>' +\n    \"<w:p><w:r><w:t>Tror jeg skal til at have noget front-end p\\u00e5, kan v\\u00e6re sv\\u00e6rt at se hvad der skal til hvis ikke jeg har et sted hvor jeg kan smide det hen.</w:t></w:r></w:p>\" +\n    \"<w:p><w:r><w:t>M\\u00e5ske starte med at lave et sted hvor man logger ind</w:t></w:r></w:p>\";\n  fikP.getRange().insertOoxml(wrapOoxml(newXml), Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Applies the diary-TODO edit via the Word COM object model:\n#  1) Merge the \"Jeg\" + \" har \" runs (and drop the proofErr spell-check\n#     markers around \"Jeg\") into a single run \"Jeg har \" in the\n#     \"...haft problemer med at f\\u00e5 lavet...\" paragraph.\n#  2) In the 17/5 TODO list: \"repository\" -> \"database CRUD for user\",\n#     and \"REST api\" -> \"REST \" + \"user API\" (two runs).\n#  3) Append a \". \" run to the \"Fik endelig hul igennem...\" paragraph and\n#     insert two new paragraphs after it (before the next, empty one).\n\n$d = $word.ActiveDocument\n\nfunction New-PkgOoxml([string]$bodyXml) {\n    return '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n            '<pkg:xmlData>' +\n                '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n                    '<w:body>' + $bodyXml + '</w:body>' +\n                '</w:document>' +\n            '</pkg:xmlData>' +\n        '</pkg:part>' +\n    '</pkg:package>'\n}\n\n$CR = [char]13\n\n# ---- Change 1: \"Jeg\" + \" har \" -> single run \"Jeg har \", proofErr removed ----\n$jegHarIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text.TrimEnd($CR)\n    if ($t.IndexOf(\"haft problemer med at f\") -ge 0) {\n        $jegHarIdx = $i\n        break\n    }\n}\nif ($jegHarIdx -eq -1) {\n    throw \"Could not locate the 'Jeg har haft problemer...' paragraph\"\n}\n$fullText = $d.Paragraphs($jegHarIdx).Range.Text.TrimEnd($CR)\n$tailText = $fullText.Substring(8)   # strip leading \"Jeg har \"\n$p1Xml = '<w:p><w:r><w:rPr><w:lang w:val=\"en-DK\"/></w:rPr><w:t xml:space=\"preserve\">Jeg har </w:t></w:r><w:r><w:t>' + $tailText + '</w:t></w:r></w:p>'\n$d.Paragraphs($jegHarIdx).Range.InsertXML((New-PkgOoxml $p1Xml))\n\n# ---- Change 2: TODO list \"repository\" / \"REST api\" (17/5 section) ----\n$repoIdx = -1\nfor ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {\n    $t1 = $d.Paragraphs($i).Range.Text.TrimEnd($CR)\n    $t2 = $d.Paragraphs($i + 1).Range.Text.TrimEnd($CR)\n    if ($t1 -eq \"repository\" -and $t2 -eq \"REST api\") {\n        $repoIdx = $i\n        break\n    }\n}\nif ($repoIdx -eq -1) {\n    throw \"Could not locate the 'repository' / 'REST api' TODO pair\"\n}\n$restIdx = $repoIdx + 1\n\n$listPPr = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr>'\n\n$repoXml = '<w:p>' + $listPPr + '<w:r><w:t>database CRUD for user</w:t></w:r></w:p>'\n$d.Paragraphs($repoIdx).Range.InsertXML((New-PkgOoxml $repoXml))\n\n$restXml = '<w:p>' + $listPPr + '<w:r><w:t xml:space=\"preserve\">REST </w:t></w:r><w:r><w:t>user API</w:t></w:r></w:p>'\n$d.Paragraphs($restIdx).Range.InsertXML((New-PkgOoxml $restXml))\n\n# ---- Change 3: append \". \" run + two new paragraphs after \"Fik endelig...\" ----\n$fikIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs($i).Range.Text.TrimEnd($CR)\n    if ($t.IndexOf(\"Fik endelig hul igennem\") -eq 0) {\n        $fikIdx = $i\n        break\n    }\n}\nif ($fikIdx -eq -1) {\n    throw \"Could not locate the 'Fik endelig hul igennem...' paragraph\"\n}\n$fikText = $d.Paragraphs($fikIdx).Range.Text.TrimEnd($CR)\n$text2 = \"Tror jeg skal til at have noget front-end p\u00e5, kan v\u00e6re sv\u00e6rt at se hvad der skal til hvis ikke jeg har et sted hvor jeg kan smide det hen.\"\n$text3 = \"M\u00e5ske starte med at lave et sted hvor man logger ind\"\n$fikNewXml = '<w:p><w:r><w:t>' + $fikText + '</w:t></w:r><w:r><w:t xml:space=\"preserve\">. </w:t></w:r></w:p>' +\n             '<w:p><w:r><w:t>' + $text2 + '</w:t></w:r></w:p>' +\n             '<w:p><w:r><w:t>' + $text3 + '</w:t></w:r></w:p>'\n$d.Paragraphs($fikIdx).Range.InsertXML((New-PkgOoxml $fikNewXml))\n"}
